$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.688.13'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').Value = '2.524.11'
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.09'
$ws.Range('E5').Value = '  +3.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.02'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.576'
$ws.Range('E7').Value = '  -1.39%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.534'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.94'
$ws.Range('E10').Value = '  -0.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0812'
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.57'
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('E13').Value = '  -3.42%  '
$ws.Range('D14').Value = '2.911.74'
$ws.Range('E14').Value = '  +0.73%  '
$ws.Range('D15').Value = '2.536.10'
$ws.Range('E15').Value = '  +0.68%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.30'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.854'
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D18').Value = '42.769.29'
$ws.Range('E18').Value = '  +1.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.92'
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.72'
$ws.Range('E20').Value = '  +4.75%  '
$ws.Range('D21').Value = '0.0₃0961'
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.73'
$ws.Range('E22').Value = '  -2.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.61'
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('E24').Value = '  +1.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.09'
$ws.Range('E25').Value = '  +2.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.71'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').Value = '  +4.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '41.70'
$ws.Range('E29').Value = '  +11.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.32'
$ws.Range('E30').Value = '  +1.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.97'
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '158.17'
$ws.Range('E32').Value = '  +2.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.15'
$ws.Range('E33').Value = '  +4.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.30'
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.68'
$ws.Range('E35').Value = '  +3.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.31'
$ws.Range('E36').Value = '  +1.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0781'
$ws.Range('E37').Value = '  -0.35%  '
$ws.Range('E38').Value = '  -2.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.119'
$ws.Range('E39').Value = '  -0.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.62'
$ws.Range('E40').Value = '  -0.63%  '
$ws.Range('E41').Value = '  +16.45%  '
$ws.Range('E42').Value = '  +1.63%  '
$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.34'
$ws.Range('E43').Value = '  -1.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.80'
$ws.Range('E44').Value = '  -1.50%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('D46').Value = '2.037.62'
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '84.59'
$ws.Range('E47').Value = '  +0.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.93'
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '75.77'
$ws.Range('E49').Value = '  +4.21%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '105.17'
$ws.Range('E50').Value = '  +4.22%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.763.75'
$ws.Range('E51').Value = '  +0.58%  '
